$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "2+14="
$t.Cell(1,2).Range.Text = "68-44="
$t.Cell(1,3).Range.Text = "61+28="
$t.Cell(1,4).Range.Text = "22+20="
$t.Cell(1,5).Range.Text = "60+26="
$t.Cell(2,1).Range.Text = "65+6="
$t.Cell(2,2).Range.Text = "15+23="
$t.Cell(2,3).Range.Text = "9-3="
$t.Cell(2,4).Range.Text = "31+12="
$t.Cell(2,5).Range.Text = "18+25="
$t.Cell(3,1).Range.Text = "26+34="
$t.Cell(3,2).Range.Text = "94-38="
$t.Cell(3,3).Range.Text = "34+21="
$t.Cell(3,4).Range.Text = "8+56="
$t.Cell(3,5).Range.Text = "76-31="
$t.Cell(4,1).Range.Text = "98-30="
$t.Cell(4,2).Range.Text = "65-48="
$t.Cell(4,3).Range.Text = "18+67="
$t.Cell(4,4).Range.Text = "92-51="
$t.Cell(4,5).Range.Text = "70-3="
$t.Cell(5,1).Range.Text = "44+20="
$t.Cell(5,2).Range.Text = "13-2="
$t.Cell(5,3).Range.Text = "66-28="
$t.Cell(5,4).Range.Text = "37-0="
$t.Cell(5,5).Range.Text = "53-4="
$t.Cell(6,1).Range.Text = "39+25="
$t.Cell(6,2).Range.Text = "14+22="
$t.Cell(6,3).Range.Text = "70-43="
$t.Cell(6,4).Range.Text = "97-24="
$t.Cell(6,5).Range.Text = "18+55="
$t.Cell(7,1).Range.Text = "79-44="
$t.Cell(7,2).Range.Text = "66-25="
$t.Cell(7,3).Range.Text = "6+81="
$t.Cell(7,4).Range.Text = "70-40="
$t.Cell(7,5).Range.Text = "10+70="
$t.Cell(8,1).Range.Text = "36+56="
$t.Cell(8,2).Range.Text = "65-44="
$t.Cell(8,3).Range.Text = "64-56="
$t.Cell(8,4).Range.Text = "9+79="
$t.Cell(8,5).Range.Text = "56+29="
$t.Cell(9,1).Range.Text = "27+18="
$t.Cell(9,2).Range.Text = "58+14="
$t.Cell(9,3).Range.Text = "58+40="
$t.Cell(9,4).Range.Text = "39-34="
$t.Cell(9,5).Range.Text = "45+52="
$t.Cell(10,1).Range.Text = "4+7="
$t.Cell(10,2).Range.Text = "55-23="
$t.Cell(10,3).Range.Text = "14+31="
$t.Cell(10,4).Range.Text = "54+22="
$t.Cell(10,5).Range.Text = "62+16="
$t.Cell(11,1).Range.Text = "68+12="
$t.Cell(11,2).Range.Text = "2+76="
$t.Cell(11,3).Range.Text = "82-2="
$t.Cell(11,4).Range.Text = "66+18="
$t.Cell(11,5).Range.Text = "70+29="
$t.Cell(12,1).Range.Text = "26+12="
$t.Cell(12,2).Range.Text = "87-58="
$t.Cell(12,3).Range.Text = "68-58="
$t.Cell(12,4).Range.Text = "31+24="
$t.Cell(12,5).Range.Text = "5+19="
$t.Cell(13,1).Range.Text = "51-28="
$t.Cell(13,2).Range.Text = "31+49="
$t.Cell(13,3).Range.Text = "58-4="
$t.Cell(13,4).Range.Text = "36+20="
$t.Cell(13,5).Range.Text = "21+59="
$t.Cell(14,1).Range.Text = "38+7="
$t.Cell(14,2).Range.Text = "61-18="
$t.Cell(14,3).Range.Text = "79-4="
$t.Cell(14,4).Range.Text = "37+47="
$t.Cell(14,5).Range.Text = "7+37="
$t.Cell(15,1).Range.Text = "67-63="
$t.Cell(15,2).Range.Text = "93+4="
$t.Cell(15,3).Range.Text = "7+4="
$t.Cell(15,4).Range.Text = "53-51="
$t.Cell(15,5).Range.Text = "0+38="
$t.Cell(16,1).Range.Text = "9-1="
$t.Cell(16,2).Range.Text = "99-93="
$t.Cell(16,3).Range.Text = "63-34="
$t.Cell(16,4).Range.Text = "74-16="
$t.Cell(16,5).Range.Text = "18+27="
$t.Cell(17,1).Range.Text = "39+59="
$t.Cell(17,2).Range.Text = "37+25="
$t.Cell(17,3).Range.Text = "40+23="
$t.Cell(17,4).Range.Text = "53-2="
$t.Cell(17,5).Range.Text = "62-29="
$t.Cell(18,1).Range.Text = "73-8="
$t.Cell(18,2).Range.Text = "60-25="
$t.Cell(18,3).Range.Text = "0+1="
$t.Cell(18,4).Range.Text = "64+17="
$t.Cell(18,5).Range.Text = "10+59="
$t.Cell(19,1).Range.Text = "57-31="
$t.Cell(19,2).Range.Text = "38-0="
$t.Cell(19,3).Range.Text = "42+11="
$t.Cell(19,4).Range.Text = "52-28="
$t.Cell(19,5).Range.Text = "2+63="
$t.Cell(20,1).Range.Text = "77+21="
$t.Cell(20,2).Range.Text = "27-17="
$t.Cell(20,3).Range.Text = "76-45="
$t.Cell(20,4).Range.Text = "76-38="
$t.Cell(20,5).Range.Text = "77-41="
